# Applies the update described by the diff:
#  - List1!A20: 45499.625 -> 45499.25
#  - List1!A21: 45500.614583333336 -> 45500.25
#  - Mereni aktivity row 19: fill in previously empty measurement row
#  - Selection (active cell) changes on both sheets (cosmetic)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("List1")
$ws2 = $wb.Worksheets.Item("Měření aktivity")

# --- Sheet List1: update the two date/time source values ---
$ws1.Range("A20").Value = 45499.25
$ws1.Range("A21").Value = 45500.25

# --- Sheet "Měření aktivity": fill in row 19 with the new measurement data ---
$ws2.Range("A19").Value = 45499.251388888886
$ws2.Range("B19").Value = 0.004
$ws2.Range("C19").Value = 0.005
$ws2.Range("D19").Value = 0.004
$ws2.Range("E19").Value = 0.003
$ws2.Range("F19").Value = 0.004
$ws2.Range("G19").Value = 152.4
$ws2.Range("H19").Value = 152.3
$ws2.Range("I19").Value = 152.3
$ws2.Range("J19").Value = 152.3
$ws2.Range("K19").Value = 152.6
$ws2.Range("L19").Value = 152.7
$ws2.Range("M19").Value = 152.6
$ws2.Range("N19").Value = 152.6
$ws2.Range("O19").Value = 152.6
$ws2.Range("P19").Value = 152.6

$excel.CalculateFullRebuild() | Out-Null

# --- Selection changes (cosmetic, reflects where the user clicked last) ---
$ws1.Activate() | Out-Null
$ws1.Range("A21").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("H20").Select() | Out-Null
